$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.387.67"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.573.53"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'290.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.3758"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("D8").Value = "'50.17"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "'0.3417"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'1.165"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'0.07675"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "'21.36"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "'5.989"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'6.930"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "1.573.70"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'90.37"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'0.06717"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").Value = "'6.242"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'0.5272"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "22.394.07"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'2.392"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "'2.772"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.44%  "
$ws.Range("D28").Value = "'20.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").Value = "'144.63"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'5.085"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").Value = "'126.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("D32").Value = "1.747.83"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "'1.025"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.64%  "
$ws.Range("D34").Value = "'6.251"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'10.07"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("D37").Value = "'0.08518"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'0.02550"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "'0.06549"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "'5.515"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'1.296"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").Value = "'11.66"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'0.6443"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'14.07"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").Value = "'0.6016"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'3.782"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'1.302"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +11.47%  "
$ws.Range("D50").Value = "'2.097"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'125.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.30%  "
